$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the new rows are treated as plain text (dates/times must stay as
# literal strings, not be auto-converted to date/time serial numbers).
$target = $ws.Range("A147:I151")
$target.NumberFormat = "@"

$data = @(
    @("2024-05-22", "09:24:06", "No pone tornillo",     "-", "-", "-", "-", "09:24:08", "0:00:02"),
    @("2024-05-22", "09:24:18", "Etiquetadora",          "-", "-", "-", "-", "09:24:19", "0:00:01"),
    @("2024-05-22", "09:26:52", "Fallo en paletizador",  "-", "-", "-", "-", "09:26:53", "0:00:01"),
    @("2024-05-22", "09:37:40", "Fallo en elevador",     "-", "-", "-", "-", "09:37:46", "0:00:06"),
    @("2024-05-22", "09:37:44", "No coge placa",         "-", "-", "-", "-", "09:37:48", "0:00:04")
)

$startRow = 147
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
